$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    'PEOPLES R CHINA',
    'USA',
    'INDIA',
    'GERMANY',
    'SOUTH KOREA',
    'JAPAN',
    'ENGLAND',
    'NETHERLANDS',
    'IRAN',
    'CANADA',
    'FRANCE',
    'ITALY',
    'AUSTRALIA',
    'BELGIUM',
    'SPAIN',
    'SINGAPORE',
    'RUSSIA',
    'TAIWAN',
    'SAUDI ARABIA',
    'TURKEY',
    'SWITZERLAND',
    'BRAZIL',
    'MALAYSIA',
    'POLAND',
    'PORTUGAL',
    'EGYPT',
    'SWEDEN',
    'PAKISTAN',
    'THAILAND',
    'ISRAEL',
    'FINLAND',
    'GREECE',
    'TURKIYE',
    'CZECH REPUBLIC',
    'DENMARK',
    'ROMANIA',
    'AUSTRIA',
    'SCOTLAND',
    'WALES',
    'CHILE',
    'QATAR',
    'VIETNAM',
    'NORWAY',
    'INDONESIA',
    'IRELAND',
    'U ARAB EMIRATES',
    'IRAQ',
    'ARGENTINA',
    'NEW ZEALAND',
    'NIGERIA',
    'SOUTH AFRICA',
    'MEXICO',
    'HUNGARY',
    'BANGLADESH',
    'SERBIA',
    'SLOVENIA',
    'UKRAINE',
    'LITHUANIA',
    'CROATIA',
    'LUXEMBOURG',
    'ETHIOPIA',
    'SLOVAKIA',
    'CYPRUS',
    'KAZAKHSTAN',
    'MOROCCO',
    'BELARUS',
    'BULGARIA',
    'JORDAN',
    'LATVIA',
    'NORTH IRELAND',
    'ALGERIA',
    'COLOMBIA',
    'PHILIPPINES',
    'TUNISIA',
    'ESTONIA',
    'UZBEKISTAN',
    'ARMENIA',
    'BOSNIA HERCEG',
    'PERU',
    'KUWAIT',
    'LIBYA',
    'AZERBAIJAN',
    'BRUNEI',
    'ECUADOR',
    'MALTA',
    'OMAN',
    'SUDAN',
    'BAHRAIN',
    'COSTA RICA',
    'KENYA',
    'KOSOVO',
    'LEBANON',
    'NORTH MACEDONIA',
    'SRI LANKA',
    'VENEZUELA',
    'YEMEN',
    'BOTSWANA',
    'GEORGIA',
    'GHANA',
    'ICELAND',
    'NEPAL',
    'PALESTINE',
    'TANZANIA',
    'ALBANIA',
    'BAHAMAS',
    'CAMBODIA',
    'COTE IVOIRE',
    'CUBA',
    'MALAWI',
    'MAURITIUS',
    'MOLDOVA',
    'MONGOLIA',
    'PARAGUAY',
    'SYRIA',
    'UGANDA',
    'URUGUAY',
    'ZAMBIA'
)

$counts = @(
    9762,
    3072,
    1027,
    979,
    888,
    739,
    707,
    562,
    514,
    511,
    497,
    463,
    455,
    382,
    382,
    288,
    233,
    229,
    202,
    192,
    191,
    186,
    169,
    158,
    144,
    138,
    135,
    128,
    97,
    94,
    85,
    81,
    81,
    80,
    78,
    73,
    72,
    68,
    65,
    60,
    60,
    60,
    48,
    46,
    46,
    46,
    44,
    43,
    41,
    38,
    34,
    31,
    30,
    26,
    26,
    25,
    24,
    23,
    18,
    17,
    16,
    16,
    15,
    15,
    13,
    12,
    12,
    11,
    11,
    11,
    10,
    10,
    8,
    8,
    7,
    7,
    6,
    6,
    6,
    5,
    5,
    4,
    4,
    4,
    4,
    4,
    4,
    3,
    3,
    3,
    3,
    3,
    3,
    3,
    3,
    3,
    2,
    2,
    2,
    2,
    2,
    2,
    2,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $counts[$i]
}
